# Weekly driver report update for 2025-04-28
# Refresh both the "Bad Drivers" and "Good Drivers" tables on the
# "Driver Summary" sheet with the latest roaming-impact numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Table 1 "Bad Drivers" (rows 3-11) - re-sorted / refreshed rows
# ---------------------------------------------------------------
$badDrivers = @(
    @{ Row = 3;  Name = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.24.1"; B = 1; C = 726;  D = 91.2 },
    @{ Row = 4;  Name = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.120.0.3";       B = 2; C = 990;  D = 92 },
    @{ Row = 5;  Name = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.20.1.1";        B = 1; C = 556;  D = 94.1 },
    @{ Row = 6;  Name = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.100.0.3";       B = 2; C = 571;  D = 94.2 },
    @{ Row = 7;  Name = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.40.0.7";        B = 1; C = 381;  D = 94.5 },
    @{ Row = 8;  Name = "Intel(R) Dual Band Wireless-AC 7265 - 19.51.12.3"; B = 3; C = 1260; D = 95.7 },
    @{ Row = 9;  Name = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.17.1"; B = 1; C = 108;  D = 96.9 },
    @{ Row = 10; Name = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.3.3";  B = 1; C = 81;   D = 97.4 },
    @{ Row = 11; Name = "Intel(R) Dual Band Wireless-AC 8265 - 20.50.3.3";  B = 8; C = 705;  D = 98.3 }
)

foreach ($row in $badDrivers) {
    $r = $row.Row
    $ws.Range("A$r").Value = $row.Name
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
}

# Totals row for table 1
$ws.Range("B12").Value = 20
$ws.Range("C12").Value = 5378

# ---------------------------------------------------------------
# Table 2 "Good Drivers" (rows 20-41) gains a brand-new entry.
# Insert a blank row at 22 (shifts existing rows 22-46 down to
# 23-47, copying formatting from the row above, same as Excel's
# native "Insert Copied Cells" / "Insert Sheet Rows" behaviour).
# ---------------------------------------------------------------
$ws.Rows("22:22").Insert()

# New row 22 - a driver version with no recorded "Driver Vintage" date
$ws.Range("A22").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B22").Value = 11128
$ws.Range("D22").Value = 100

# Refreshed sample counts for rows that shifted down (text/dates unchanged)
$ws.Range("B23").Value = 486214
$ws.Range("B25").Value = 79953
$ws.Range("B26").Value = 35355
$ws.Range("B27").Value = 172690
$ws.Range("B29").Value = 65425
$ws.Range("B30").Value = 236471
$ws.Range("B31").Value = 41618
$ws.Range("B32").Value = 117653
$ws.Range("B35").Value = 154175
$ws.Range("B40").Value = 132352

# ---------------------------------------------------------------
# Nudge the sheet's recorded used-range so it covers the newly
# appended trailing blank row together with the existing spacer
# columns through J, matching the refreshed report's extent
# (A1:J47), touching only the single far corner cell.
# ---------------------------------------------------------------
$ws.Range("J47").HorizontalAlignment = $ws.Range("J47").HorizontalAlignment
